$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Edit and add candidate dashboard": the sheet had a spare "true/false"
# flag column (AA) ahead of the last data column (AD). Remove that column
# entirely -- everything to its right (AB:AD) shifts one column left
# (new AA:AC), shrinking the used range from A1:AD3 to A1:AC3.
$ws.Columns.Item(27).Delete()

# The last column (old AD, now AC) holds a pin/OTP-like code. Row 2's
# value changes from "553321" to "5555". Force text so the numeric-looking
# string isn't reinterpreted as a number (matches the existing shared
# text cell it replaces).
$cell = $ws.Range("AC2")
$cell.NumberFormat = "@"
$cell.Value = "5555"
